# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
# Both sheets carry the same exhibition data and received the same updates.

$wb = $excel.ActiveWorkbook

$updates = @{
    2 = 1951
    3 = 587
    4 = 1291
    5 = 6463
    6 = 155
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
